# New crime data collected — refresh the weekly CompStat report
# (17th Precinct, report week / volume header text, and the weekly
# crime-complaint statistics table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (shared strings with rich-text runs) ---
# "Volume 31   Number  17" -> "...18"
$ws.Range("A8").Value = "Volume 31   Number  18"
# "Report Covering the Week  4/22/2024  Through  4/28/2024" -> new week
$ws.Range("C9").Value = "Report Covering the Week  4/29/2024  Through  5/5/2024"

# --- Simple numeric value updates (style/number-format unchanged) ---
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 100

$ws.Range("C16").Value = 2
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 29
$ws.Range("K16").Value = 3.571428571428
$ws.Range("L16").Value = 26.086956521739
$ws.Range("M16").Value = 31.818181818181
$ws.Range("N16").Value = -85.128205128205

$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 180
$ws.Range("I17").Value = 47
$ws.Range("J17").Value = 38
$ws.Range("K17").Value = 23.684210526315
$ws.Range("L17").Value = 11.904761904761
$ws.Range("M17").Value = 104.347826086957
$ws.Range("N17").Value = -17.543859649122

$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 37
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = -19.565217391304
$ws.Range("L18").Value = -43.076923076923
$ws.Range("M18").Value = -9.756097560975
$ws.Range("N18").Value = -90.703517587939

$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = -50.666666666666
$ws.Range("I19").Value = 200
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = -20
$ws.Range("L19").Value = -5.213270142180
$ws.Range("M19").Value = -15.611814345991
$ws.Range("N19").Value = -70.414201183432

$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -75
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -52.941176470588
$ws.Range("L20").Value = -38.461538461538
$ws.Range("N20").Value = -96.506550218340

$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -11.111111111111
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -26.262626262626
$ws.Range("I21").Value = 326
$ws.Range("J21").Value = 381
$ws.Range("K21").Value = -14.435695538057
$ws.Range("L21").Value = -8.683473389355
$ws.Range("M21").Value = -1.510574018126
$ws.Range("N21").Value = -79.089159717767

$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 150

$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -16.091954022988
$ws.Range("I24").Value = 330
$ws.Range("J24").Value = 324
$ws.Range("K24").Value = 1.851851851851
$ws.Range("L24").Value = -18.518518518518
$ws.Range("M24").Value = 53.488372093023

$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -23.529411764705
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = -15.714285714285
$ws.Range("I25").Value = 275
$ws.Range("J25").Value = 261
$ws.Range("K25").Value = 5.363984674329
$ws.Range("L25").Value = -18.639053254437

$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -71.428571428571
$ws.Range("F26").Value = 10
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -64.285714285714
$ws.Range("I26").Value = 71
$ws.Range("J26").Value = 93
$ws.Range("K26").Value = -23.655913978494
$ws.Range("L26").Value = -11.25
$ws.Range("M26").Value = -6.578947368421

$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 25

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 5
$ws.Range("K31").Value = 66.666666666666
$ws.Range("L31").Value = 0

# --- Cells that change numeric <-> "no data" text placeholder type ---
# Shared placeholder cells: text "0" (e.g. C14) and text "***.*" (e.g. E14)
# used throughout the table to represent zero / not-applicable values.
# Copying from an existing placeholder cell keeps the exact style (incl.
# number format/font) used for these placeholder cells.

# Row 28: count goes from 2 to 1; pct-change columns become "no data"
$ws.Range("C28").Value = 1
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))

# Row 31: D/E become "no data" placeholders; F becomes a real number
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
$ws.Range("D15").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 1

# Rows 29, 30 and 33: "no data" placeholder -> real computed percentage
$ws.Range("E15").Copy($ws.Range("N29"))
$ws.Range("N29").Value = -100
$ws.Range("E15").Copy($ws.Range("N30"))
$ws.Range("N30").Value = -100
$ws.Range("E15").Copy($ws.Range("L33"))
$ws.Range("L33").Value = -100
